# repull data, push all data, mean calculation
# Update column F (dSF) values for specific rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = -8
    3  = -2
    5  = 3
    10 = 2
    12 = 5
    13 = 4
    30 = 2
    38 = 0
    43 = -1
    44 = 1
    47 = 1
    49 = -1
    51 = 3
    52 = 2
    58 = -2
    67 = 0
    68 = -2
    69 = 1
    78 = 0
    80 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
